$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(6)

# Insert a new column before I. Old I (date) -> J, old J (legislator_name) -> K,
# old K (legislator_id) -> L. The new, empty column I becomes "category".
$ws.Columns.Item(9).Insert()

# Copy formatting (border/bold/alignment for header, plain for data) from the
# legislator_id column (L) onto the two new trailing columns M (source_file)
# and N (index) so they pick up the same per-row look.
$ws.Range("L1:L13").Copy()
$ws.Range("M1:M13").PasteSpecial(-4122)
$ws.Range("L1:L13").Copy()
$ws.Range("N1:N13").PasteSpecial(-4122)

# Header row
$ws.Range("I1").Value = "category"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# Data rows 2..13: category="normal", source_file="tmp671f1", index=A (same row id)
for ($r = 2; $r -le 13; $r++) {
    $ws.Cells.Item($r, 9).Value = "normal"
    $ws.Cells.Item($r, 13).Value = "tmp671f1"
    $idxVal = $ws.Cells.Item($r, 1).Value()
    $ws.Cells.Item($r, 14).Value = $idxVal
}
